# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N ("Late"), matching column M's width, then make this sheet the
# active one with cell S6 selected (mirrors the author switching tabs and
# inserting a column while editing the repayment schedule).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of column M (13) so the freshly inserted column N
# matches it, same as Excel does when inserting a column.
$existingWidth = $ws.Columns.Item(13).ColumnWidth

# Insert a blank column before the current column N ("Late"); everything
# from N onward shifts one column to the right (N->O, O->P, P->Q).
$ws.Columns.Item(14).Insert()

# Give the newly inserted column N the same width as column M.
$ws.Columns.Item(14).ColumnWidth = $existingWidth

# Make "Repayment schedule" the active sheet/tab with S6 selected.
$null = $ws.Activate()
$null = $ws.Range("S6").Select()
